$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-27 Sunday", "2024-10-28 Monday"),
    @("735×5=", "676×4="),
    @("319×2=", "724×4="),
    @("248×8=", "759×2="),
    @("432×5=", "721×7="),
    @("346×5=", "278×4="),
    @("357×5=", "203×5="),
    @("353×2=", "994×8="),
    @("653×4=", "718×8="),
    @("854×7=", "492×7="),
    @("855×4=", "434×7="),
    @("354×2=", "223×2="),
    @("322×8=", "888×7="),
    @("439×5=", "498×2="),
    @("736×9=", "134×6="),
    @("571×5=", "999×6="),
    @("467×9=", "898×9="),
    @("238×3=", "421×9="),
    @("210×8=", "946×2="),
    @("983×4=", "856×3="),
    @("104×8=", "219×4="),
    @("523×5=", "325×6="),
    @("649×7=", "730×4="),
    @("273×8=", "714×5="),
    @("179×2=", "986×5="),
    @("819×6=", "489×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
